$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.942.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.219.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.397"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.700"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.52%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.216.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.585"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.181"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.819.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.862.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.236.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("B19").Value = "PEPE"
$ws.Range("C19").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000236"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +88.36%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +19.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "442.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.27%  "
$ws.Range("E23").Value = "  +6.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "83.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.454.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.161"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +45.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "553.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.50%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +6.67%  "
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "175.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.781"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.70%  "
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.631"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.19%  "
